$d = $word.ActiveDocument

$d.Content.Find.Execute("analyse", $true, $false, $false, $false, $false, $true, 1, $false, "analyze", 2)

$d.Content.Find.Execute("github.", $true, $false, $false, $false, $false, $true, 1, $false,
    "zenodo (https://zenodo.org/record/7765007).", 2)
